# Simplify file upload functionality
# Append a new data row (row 43) to each of the four log sheets.

$wb = $excel.ActiveWorkbook

$sheetsData = @{
    "DE_LFT_#1" = @{
        A = 45829.43844907408
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        D = "0x01,0x68"
        E = "0x14"
        F = 380
        G = 759863127514710945038336.0
        H = 360
        I = 14
    }
    "DE_LFT_#2" = @{
        A = 45829.43844907408
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        D = "0x01,0x68"
        E = "0xe"
        F = 380
        G = 568432987514711010443264.0
        H = 360
        I = 14
    }
    "DE_PLT_#1" = @{
        A = 45829.43844907408
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x00,0x7F"
        E = "0x7"
        F = 130
        G = 568631262647113970876416.0
        H = 127
        I = 7
    }
    "DE_PLT_#2" = @{
        A = 45829.43844907408
        B = "0x00,0x82"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x00,0x7F"
        E = "0x3"
        F = 130
        G = 985046333984776009023488.0
        H = 127
        I = 3
    }
}

$sheetNames = @("DE_LFT_#1", "DE_LFT_#2", "DE_PLT_#1", "DE_PLT_#2")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $data = $sheetsData[$sheetName]
    $row = 43

    $ws.Cells.Item($row, 1).Value = $data.A
    $ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 2).Value = $data.B
    $ws.Cells.Item($row, 3).Value = $data.C
    $ws.Cells.Item($row, 4).Value = $data.D
    $ws.Cells.Item($row, 5).Value = $data.E
    $ws.Cells.Item($row, 6).Value = $data.F
    $ws.Cells.Item($row, 7).Value = $data.G
    $ws.Cells.Item($row, 8).Value = $data.H
    $ws.Cells.Item($row, 9).Value = $data.I
}
